$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enable iterative calculation with the new max-change (iterateDelta) setting.
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# Row 33: 30-Jan-2012, EKA_METALS_PATCH_0357
$ws.Range("B33").Value = 40938
$ws.Range("B33").NumberFormat = "d-mmm"
$ws.Range("I33").Value = "EKA_METALS_PATCH_0357"

# Row 34: 31-Jan-2012, EKA_METALS_PATCH_0382
$ws.Range("B34").Value = 40939
$ws.Range("B34").NumberFormat = "d-mmm"
$ws.Range("I34").Value = "EKA_METALS_PATCH_0382"

# Row 35: 9-Feb-2012 (DB Patch applied)
$ws.Range("B35").Value = 40948
$ws.Range("B35").NumberFormat = "d-mmm"

# Update the active selection to I34, matching the authored workbook state.
$ws.Range("I34").Select()
